# Update the cryptos price list (columns D = Price, E = Volume(1h)) with the
# latest scraped values. Cells whose new Price text looks like a plain
# decimal number (e.g. "22.83") are prefixed with a leading apostrophe so
# Excel stores them as text (matching the existing column formatting)
# instead of silently converting them to a numeric value.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "27.172.68"
$ws.Range("E2").Value = "  +0.47%  "
$ws.Range("D3").Value = "1.684.30"
$ws.Range("E3").Value = "  +0.10%  "
$ws.Range("E4").Value = "  +0.09%  "
$ws.Range("D5").Value = "'216.06"
$ws.Range("E5").Value = "  -0.03%  "
$ws.Range("E6").Value = "  -0.01%  "
$ws.Range("E7").Value = "  +0.11%  "
$ws.Range("D8").Value = "'22.83"
$ws.Range("E8").Value = "  +6.22%  "
$ws.Range("E9").Value = "  +2.39%  "
$ws.Range("E10").Value = "  +0.24%  "
$ws.Range("E11").Value = "  -0.01%  "
$ws.Range("D12").Value = "1.922.46"
$ws.Range("E12").Value = "  +0.11%  "
$ws.Range("D13").Value = "1.681.78"
$ws.Range("E13").Value = "  -0.29%  "
$ws.Range("E14").Value = "  +2.07%  "
$ws.Range("E15").Value = "  +4.58%  "
$ws.Range("D16").Value = "'66.92"
$ws.Range("E16").Value = "  +0.77%  "
$ws.Range("D17").Value = "27.181.93"
$ws.Range("E17").Value = "  +0.46%  "
$ws.Range("D18").Value = "'235.87"
$ws.Range("E18").Value = "  -0.27%  "
$ws.Range("D19").Value = "'7.96"
$ws.Range("E19").Value = "  -2.72%  "
$ws.Range("E20").Value = "  +0.37%  "
$ws.Range("E21").Value = "  +0.08%  "
$ws.Range("D22").Value = "'4.56"
$ws.Range("E22").Value = "  +1.85%  "
$ws.Range("D23").Value = "'9.55"
$ws.Range("E23").Value = "  +2.85%  "
$ws.Range("D24").Value = "'2.08"
$ws.Range("E24").Value = "  -2.21%  "
$ws.Range("D25").Value = "'147.07"
$ws.Range("E25").Value = "  -0.09%  "
$ws.Range("D26").Value = "'7.41"
$ws.Range("E26").Value = "  +2.24%  "
$ws.Range("D27").Value = "'16.41"
$ws.Range("E27").Value = "  -2.03%  "
$ws.Range("E28").Value = "  +0.12%  "
$ws.Range("E29").Value = "  +0.04%  "
$ws.Range("E30").Value = "  +1.04%  "
$ws.Range("E31").Value = "  +0.10%  "
$ws.Range("E32").Value = "  +0.42%  "
$ws.Range("D33").Value = "1.548.73"
$ws.Range("E33").Value = "  +1.15%  "
$ws.Range("E34").Value = "  +2.09%  "
$ws.Range("E35").Value = "  -2.17%  "
$ws.Range("D36").Value = "'0.605"
$ws.Range("E36").Value = "  +2.42%  "
$ws.Range("D37").Value = "'0.948"
$ws.Range("E37").Value = "  +3.11%  "
$ws.Range("E38").Value = "  -0.37%  "
$ws.Range("E39").Value = "  -1.29%  "
$ws.Range("E40").Value = "  +2.43%  "
$ws.Range("D41").Value = "'5.74"
$ws.Range("E41").Value = "  -0.09%  "
$ws.Range("D42").Value = "'69.06"
$ws.Range("E42").Value = "  +1.47%  "
$ws.Range("E43").Value = "  +0.12%  "
$ws.Range("E44").Value = "  -0.43%  "
$ws.Range("D45").Value = "1.829.16"
$ws.Range("E45").Value = "  +0.23%  "
$ws.Range("D46").Value = "'0.791"
$ws.Range("E46").Value = "  +1.33%  "
$ws.Range("D47").Value = "'89.72"
$ws.Range("E47").Value = "  -0.77%  "
$ws.Range("E48").Value = "  +7.81%  "
$ws.Range("D49").Value = "'1.62"
$ws.Range("E49").Value = "  +5.83%  "
$ws.Range("D50").Value = "'8.22"
$ws.Range("E50").Value = "  +4.07%  "
$ws.Range("E51").Value = "  -0.57%  "
